# Auto-generated Excel COM-interop script applying the cryptos.xlsx data refresh
# (Price column D and Volume(1h) column E updates) described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text would be mis-parsed as a number by Excel (e.g. "236.23",
# "1.000", "0.9999") need NumberFormat temporarily forced to Text ("@") so the
# literal string is preserved (no lost trailing zeros / no numeric conversion).
# The original style is captured first and restored afterward so no visible
# formatting / style-index change is left behind.
$textForceCells = @(
    "D5"
    "D6"
    "D7"
    "D8"
    "D11"
    "D12"
    "D13"
    "D14"
    "D15"
    "D16"
    "D19"
    "D20"
    "D22"
    "D23"
    "D24"
    "D25"
    "D26"
    "D28"
    "D29"
    "D31"
    "D33"
    "D34"
    "D35"
    "D36"
    "D37"
    "D38"
    "D40"
    "D41"
    "D42"
    "D43"
    "D44"
    "D45"
    "D46"
    "D48"
    "D51"
)

$origStyles = @{}
foreach ($addr in $textForceCells) {
    $origStyles[$addr] = $ws.Range($addr).Style
    $ws.Range($addr).NumberFormat = "@"
}

# --- Apply every cell update (Price column D, Volume(1h) column E) in row order.
$ws.Range("D2").Value = "26.151.56"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.750.88"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "236.23"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "0.5301"
$ws.Range("E7").Value = "  +1.66%  "
$ws.Range("D8").Value = "0.2808"
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").Value = "1.740.35"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "0.07177"
$ws.Range("E11").Value = "  +2.01%  "
$ws.Range("D12").Value = "15.46"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "0.6480"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "4.631"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "78.48"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").Value = "0.9999"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "26.039.12"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "11.76"
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").Value = "0.000006771"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("D21").Value = "1.969.81"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").Value = "4.344"
$ws.Range("E22").Value = "  +4.18%  "
$ws.Range("D23").Value = "8.732"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").Value = "5.243"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "139.14"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "1.529"
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").Value = "1.799"
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("D29").Value = "104.91"
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "3.775"
$ws.Range("E31").Value = "  +3.18%  "
$ws.Range("E32").Value = "  +5.81%  "
$ws.Range("D33").Value = "0.04638"
$ws.Range("E33").Value = "  +3.74%  "
$ws.Range("D34").Value = "2.643"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("D35").Value = "1.008"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("D36").Value = "0.6309"
$ws.Range("E36").Value = "  +3.26%  "
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").Value = "0.01620"
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("D40").Value = "0.9999"
$ws.Range("D41").Value = "101.96"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("D42").Value = "0.3925"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("D43").Value = "0.7565"
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("D44").Value = "5.085"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "0.1153"
$ws.Range("E45").Value = "  +2.93%  "
$ws.Range("D46").Value = "6.355"
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("E47").Value = "  -2.42%  "
$ws.Range("D48").Value = "54.67"
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").Value = "7.630"
$ws.Range("E51").Value = "  -0.16%  "

# --- restore original styles on the text-forced cells so no stray style index is left
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = $origStyles[$addr]
}
